$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the three new rows of data -------------------------------------
$ws.Range("A5").Value2 = 4
$ws.Range("B5").Value2 = 44606
$ws.Range("C5").Value2 = "Харпер Ли"
$ws.Range("D5").Value2 = "Убить пересмешника"
$ws.Range("E5").Value2 = "Стр 154."

$ws.Range("A6").Value2 = 5
$ws.Range("B6").Value2 = 44641
$ws.Range("C6").Value2 = "Джейн Остен"
$ws.Range("D6").Value2 = " Гордость и предубеждение"
$ws.Range("E6").Value2 = "Стр 123."

$ws.Range("A7").Value2 = 6
$ws.Range("B7").Value2 = 44665
$ws.Range("C7").Value2 = "Анна Франк"
$ws.Range("D7").Value2 = "Дневник Анны Франк"
$ws.Range("E7").Value2 = "Стр 243."

# --- 2. Make sure the date column isn't carrying a "built-in" number format
#        before we touch borders (engine quirk: touching style on a cell
#        that already has numFmtId 14 re-derives it as a custom format) ----
$dateCol = $ws.Range("B2:B7")
$dateCol.NumberFormat = "General"

# --- 3. Table borders: thin grid inside, medium box around the outside ----
$full = $ws.Range("A1:E7")
$full.Borders(11).Weight = 2
$full.Borders(12).Weight = 2
$full.Borders(7).Weight = -4138
$full.Borders(8).Weight = -4138
$full.Borders(9).Weight = -4138
$full.Borders(10).Weight = -4138

# --- 4. Header formatting ---------------------------------------------------
$header = $ws.Range("A1:E1")
$header.Interior.Color = 15773696
$header.Font.ThemeColor = 2
$header.HorizontalAlignment = -4131
$header.VerticalAlignment = -4108

# --- 5. Re-apply the date number format now that borders are settled ------
$dateCol.NumberFormat = "mm-dd-yy"

Write-Host "done"
